$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '28.470.60'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.567.89'
$ws.Range("E3").Value = '  -2.11%  '
$ws.Range("E4").Value = '  +0.11%  '
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '212.05'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  -1.29%  '
$ws.Range("E6").Value = '  -1.15%  '
$ws.Range("E7").Value = '  +0.12%  '
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '46.26'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  +4.88%  '
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '24.01'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  +0.06%  '
$ws.Range("E10").Value = '  -1.77%  '
$ws.Range("E11").Value = '  -1.49%  '
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = '1.791.70'
$ws.Range("D14").Value = '1.567.96'
$ws.Range("E14").Value = '  -2.32%  '
$ws.Range("E15").Value = '  -2.67%  '
$ws.Range("E16").Value = '  -2.78%  '
$ws.Range("D17").Value = '28.485.79'
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '62.28'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -1.63%  '
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '227.43'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -1.87%  '
$ws.Range("E20").Value = '  -2.62%  '
$ws.Range("E21").Value = '  -2.91%  '
$ws.Range("E22").Value = '  +0.13%  '
$ws.Range("E23").Value = '  -5.75%  '
$ws.Range("E24").Value = '  -3.13%  '
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '2.11'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +7.39%  '
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = '150.93'
$cell.Style = "Normal"
$ws.Range("E26").Value = '  -1.14%  '
$ws.Range("E27").Value = '  -2.13%  '
$ws.Range("E28").Value = '  -2.62%  '
$ws.Range("E29").Value = '  -4.12%  '
$ws.Range("E30").Value = '  +0.10%  '
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '0.0469'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  -1.43%  '
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '3.22'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -1.22%  '
$ws.Range("E34").Value = '  -2.96%  '
$ws.Range("D35").Value = '1.395.14'
$ws.Range("E35").Value = '  -1.59%  '
$ws.Range("E36").Value = '  -0.75%  '
$ws.Range("E37").Value = '  -3.37%  '
$ws.Range("E38").Value = '  +1.50%  '
$ws.Range("E39").Value = '  +2.62%  '
$ws.Range("E40").Value = '  -0.71%  '
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '0.534'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  -2.20%  '
$ws.Range("E42").Value = '  +0.15%  '
$ws.Range("E43").Value = '  -4.32%  '
$ws.Range("E44").Value = '  +1.81%  '
$ws.Range("E45").Value = '  -4.23%  '
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.974'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.96%  '
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '62.76'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  -3.38%  '
$ws.Range("D48").Value = '1.703.43'
$ws.Range("E48").Value = '  -2.20%  '
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("E50").Value = '  -0.30%  '
$ws.Range("D51").Value = '0.0₆0102'
$ws.Range("E51").Value = '  -4.80%  '
